$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.163.05"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.849.52"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'235.09"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "'0.4697"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").Value = "'0.2886"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").Value = "'0.06525"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").Value = "'21.68"
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D11").Value = "'0.07950"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "'97.35"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "1.855.58"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "'5.079"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "'0.6734"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'266.54"
$ws.Range("E16").Value = "  -4.76%  "
$ws.Range("D17").Value = "30.137.83"
$ws.Range("D18").Value = "'13.57"
$ws.Range("E18").Value = "  +7.09%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007560"
$ws.Range("E19").Value = "  +4.16%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "2.099.09"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'5.191"
$ws.Range("E23").Value = "  -5.30%  "
$ws.Range("D24").Value = "'6.123"
$ws.Range("E24").Value = "  -0.26%  "
$ws.Range("D25").Value = "'166.61"
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").Value = "'9.142"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'18.79"
$ws.Range("E27").Value = "  -1.58%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +1.28%  "
$ws.Range("D30").Value = "'0.09833"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").Value = "'1.464"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "'4.264"
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").Value = "'3.991"
$ws.Range("E33").Value = "  -2.51%  "
$ws.Range("D34").Value = "'0.04679"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").Value = "'0.6963"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").Value = "'0.01859"
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").Value = "'2.596"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").Value = "'73.02"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").Value = "'0.8357"
$ws.Range("E44").Value = "  -1.07%  "
$ws.Range("D45").Value = "'103.05"
$ws.Range("D46").Value = "'0.4114"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").Value = "'938.16"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").Value = "'9.159"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").Value = "'6.977"
$ws.Range("E49").Value = "  -2.62%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "'0.05652"
